# Historias de usuario actualizadas
# Update the "Requerimientos" sheet: renumber the RE_0xx requirement codes
# in column B for rows 11-21 so they form a continuous sequence
# (RE_008 .. RE_018), fixing the previously-blank "RE_" placeholder rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Requerimientos")
$ws.Activate()

$ws.Range("B11").Value = "RE_008"
$ws.Range("B12").Value = "RE_009"
$ws.Range("B13").Value = "RE_010"
$ws.Range("B14").Value = "RE_011"
$ws.Range("B15").Value = "RE_012"
$ws.Range("B16").Value = "RE_013"
$ws.Range("B17").Value = "RE_014"
$ws.Range("B18").Value = "RE_015"
$ws.Range("B19").Value = "RE_016"
$ws.Range("B20").Value = "RE_017"
$ws.Range("B21").Value = "RE_018"

# Update the selection to match the saved cursor position
$ws.Range("D23").Select()

# Column B width now auto-fits its (shorter) contents
$ws.Columns.Item(2).AutoFit()
